$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (row 1 title)
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 02:14"

# Apply updated COVID-19 country statistics

# Row 4
$ws.Range("B4").Value = 2890409
$ws.Range("C4").Value = 54725
$ws.Range("D4").Value = 1210792
$ws.Range("E4").Value = 1547516
$ws.Range("G4").Value = 616
$ws.Range("H4").Value = 132101

# Row 5
$ws.Range("B5").Value = 1543341
$ws.Range("C5").Value = 41988
$ws.Range("D5").Value = 945915
$ws.Range("E5").Value = 534172
$ws.Range("G5").Value = 1264
$ws.Range("H5").Value = 63254

# Row 28
$ws.Range("B28").Value = 71419
$ws.Range("C28").Value = 147
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 5420

# Row 73
$ws.Range("B73").Value = 8921
$ws.Range("C73").Value = 19
$ws.Range("E73").Value = 532

# Row 97
$ws.Range("B97").Value = 3918
$ws.Range("C97").Value = 130
$ws.Range("D97").Value = 859
$ws.Range("E97").Value = 3011

# Row 125
$ws.Range("B125").Value = 1498
$ws.Range("C125").Value = 156
$ws.Range("E125").Value = 1165

# Row 135
$ws.Range("D135").Value = 839
$ws.Range("E135").Value = 141

# Row 137
$ws.Range("B137").Value = 952
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 837
$ws.Range("E137").Value = 87

# Row 140
$ws.Range("A140").Value = "Libia"
$ws.Range("B140").Value = 918
$ws.Range("C140").Value = 27
$ws.Range("D140").Value = 230
$ws.Range("E140").Value = 661
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 27

# Row 141
$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 911
$ws.Range("C141").Value = 9
$ws.Range("D141").Value = 849
$ws.Range("E141").Value = 62
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 0

# Row 142
$ws.Range("A142").Value = "Suazilandia"
$ws.Range("B142").Value = 909
$ws.Range("C142").Value = 36
$ws.Range("D142").Value = 515
$ws.Range("E142").Value = 381
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 13

# Row 172
$ws.Range("A172").Value = "Burundi"
$ws.Range("B172").Value = 191
$ws.Range("C172").Value = 21
$ws.Range("D172").Value = 118
$ws.Range("E172").Value = 72
$ws.Range("H172").Value = 1

# Row 173
$ws.Range("A173").Value = "Islas Feroe"
$ws.Range("B173").Value = 187
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 187
$ws.Range("E173").Value = 0
$ws.Range("H173").Value = 0

# Row 174
$ws.Range("A174").Value = "Guadalupe"
$ws.Range("B174").Value = 184
$ws.Range("C174").Value = 2
$ws.Range("D174").Value = 157
$ws.Range("E174").Value = 13
$ws.Range("H174").Value = 14

# Row 175
$ws.Range("A175").Value = "Gibraltar"
$ws.Range("B175").Value = 178
$ws.Range("D175").Value = 176
$ws.Range("E175").Value = 2
$ws.Range("H175").Value = 0

# Row 188
$ws.Range("B188").Value = 68
$ws.Range("E188").Value = 42
